# Apply updated crypto price/volume figures to sheet1 (per commit diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '48.298.43'
$ws.Range("E2").Value = '  +1.49%  '
$ws.Range("D3").Value = '2.507.74'
$ws.Range("E3").Value = '  +0.52%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.71'
$ws.Range("E5").Value = '  +0.00%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.21'
$ws.Range("E6").Value = '  -0.70%  '
$ws.Range("E7").Value = '  +1.22%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  -0.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.80'
$ws.Range("E10").Value = '  +0.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.27'
$ws.Range("E11").Value = '  +8.86%  '
$ws.Range("E12").Value = '  +0.99%  '
$ws.Range("E13").Value = '  +0.03%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.19'
$ws.Range("E14").Value = '  -0.20%  '
$ws.Range("D15").Value = '2.900.12'
$ws.Range("E15").Value = '  +0.51%  '
$ws.Range("D16").Value = '2.507.02'
$ws.Range("E16").Value = '  +0.44%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.844'
$ws.Range("E17").Value = '  -0.53%  '
$ws.Range("D18").Value = '48.148.84'
$ws.Range("E18").Value = '  +1.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.10'
$ws.Range("E19").Value = '  -2.05%  '
$ws.Range("E20").Value = '  +2.50%  '
$ws.Range("E21").Value = '  +0.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.78'
$ws.Range("E22").Value = '  +1.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '281.35'
$ws.Range("E23").Value = '  +13.94%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.37'
$ws.Range("E24").Value = '  +2.40%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.55'
$ws.Range("E25").Value = '  -0.21%  '
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.26'
$ws.Range("E28").Value = '  +8.66%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.79'
$ws.Range("E29").Value = '  -2.05%  '
$ws.Range("E30").Value = '  +1.12%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.35'
$ws.Range("E31").Value = '  +1.43%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.49'
$ws.Range("E32").Value = '  -0.87%  '
$ws.Range("E33").Value = '  -3.37%  '
$ws.Range("E34").Value = '  +0.25%  '
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0784'
$ws.Range("E36").Value = '  -0.59%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.95'
$ws.Range("E37").Value = '  -0.42%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.65'
$ws.Range("E38").Value = '  -1.75%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.91'
$ws.Range("E39").Value = '  -1.06%  '
$ws.Range("E40").Value = '  -0.18%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '122.14'
$ws.Range("E41").Value = '  +2.54%  '
$ws.Range("E42").Value = '  -0.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.54'
$ws.Range("E43").Value = '  -4.24%  '
$ws.Range("E44").Value = '  +1.92%  '
$ws.Range("D45").Value = '2.018.66'
$ws.Range("E45").Value = '  +0.98%  '
$ws.Range("E46").Value = '  +4.50%  '
$ws.Range("E47").Value = '  +3.05%  '
$ws.Range("E48").Value = '  -2.44%  '
$ws.Range("E49").Value = '  -0.28%  '
$ws.Range("E50").Value = '  -0.95%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '80.70'
$ws.Range("E51").Value = '  +3.88%  '
